# Refresh the crypto price/volume table (cols B-E, rows 2-51) with the
# latest scrape. Two coins (row 20/21 and row 25/26) swapped rank order,
# so their Coin/Link/Price/Volume cells are rewritten wholesale; every
# other changed row only updates Price (D) and/or Volume(1h) (E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values are plain decimals (e.g. "4.40", "1.00") that Excel
# would otherwise auto-convert to numbers (dropping the trailing zero /
# changing the stored type). Column D holds formatted price text, so force
# those specific cells to Text before writing the new value.
$textPriceCells = @('D5', 'D6', 'D8', 'D12', 'D17', 'D19', 'D20', 'D21', 'D23', 'D25', 'D26', 'D27', 'D31', 'D32', 'D34', 'D36', 'D39', 'D41', 'D43', 'D46', 'D47', 'D48', 'D50')
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = @{
    'D2' = '72.400.94'
    'E2' = '  +1.97%  '
    'D3' = '2.645.06'
    'E3' = '  +0.94%  '
    'E4' = '  +0.02%  '
    'D5' = '603.50'
    'E5' = '  -0.22%  '
    'D6' = '180.76'
    'E6' = '  -0.15%  '
    'E7' = '  -0.05%  '
    'D8' = '0.525'
    'E8' = '  +0.26%  '
    'E9' = '  +7.67%  '
    'D10' = '2.643.60'
    'E10' = '  +0.93%  '
    'E11' = '  +1.59%  '
    'D12' = '0.357'
    'E12' = '  +2.72%  '
    'E13' = '  +0.56%  '
    'E14' = '  +4.67%  '
    'D15' = '3.127.96'
    'E15' = '  +1.89%  '
    'D16' = '72.260.06'
    'E16' = '  +1.78%  '
    'D17' = '26.61'
    'E17' = '  -0.38%  '
    'D18' = '2.638.51'
    'E18' = '  +0.58%  '
    'D19' = '11.99'
    'E19' = '  +4.66%  '
    'B20' = 'BitcoinCash'
    'C20' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D20' = '380.04'
    'E20' = '  -0.48%  '
    'B21' = 'Uniswap'
    'C21' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D21' = '7.95'
    'E21' = '  +0.20%  '
    'E22' = '  +0.47%  '
    'D23' = '2.07'
    'E23' = '  +10.62%  '
    'E24' = '  +1.54%  '
    'B25' = 'NEARProtocol'
    'C25' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D25' = '4.40'
    'E25' = '  -0.95%  '
    'B26' = 'Dai'
    'C26' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D26' = '1.00'
    'E26' = '  -0.11%  '
    'D27' = '10.02'
    'E27' = '  +3.52%  '
    'D28' = '2.780.92'
    'E28' = '  +1.09%  '
    'E29' = '  +0.16%  '
    'D30' = '0.0₃0960'
    'E30' = '  +1.30%  '
    'D31' = '525.15'
    'E31' = '  -0.29%  '
    'D32' = '8.16'
    'E32' = '  +1.57%  '
    'E33' = '  -0.53%  '
    'D34' = '1.83'
    'E34' = '  -0.41%  '
    'E35' = '  -0.02%  '
    'D36' = '164.58'
    'E36' = '  -0.06%  '
    'E37' = '  +1.00%  '
    'E38' = '  -6.06%  '
    'D39' = '19.09'
    'E39' = '  +0.80%  '
    'E40' = '  +2.25%  '
    'D41' = '1.86'
    'E41' = '  -3.04%  '
    'E42' = '  +2.82%  '
    'D43' = '5.11'
    'E43' = '  +1.37%  '
    'E44' = '  +0.03%  '
    'E45' = '  +0.86%  '
    'D46' = '39.35'
    'E46' = '  -2.04%  '
    'D47' = '151.48'
    'E47' = '  -1.46%  '
    'D48' = '3.73'
    'E48' = '  +2.04%  '
    'E49' = '  +2.29%  '
    'D50' = '1.71'
    'E50' = '  +2.60%  '
    'D51' = '0.0₆0261'
    'E51' = '  -3.54%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
